$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 37 (record #36) - email written first so the shared-string table
# grows in the same order as the source edit.
$ws.Range("A37").Value = 36
$ws.Range("D37").Value = "anitasalet2203@gmail.com"

# New header for column K ("Raro"), matching the style used by the other
# header cells on row 1 (e.g. J1).
$ws.Range("K1").Value = "Raro"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B37").Value = "Ana Salet Hidalgo Flores"

# New D2 value (added alongside existing row 2 data)
$ws.Range("D2").Value = "Joselyn Dayanna Allauca Chasipanta"

$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1
$ws.Range("J37").Value = 1
$ws.Range("K37").Value = 1
$ws.Range("K37").Font.Underline = $true

$ws.Range("K37").Select()
